$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''43.769.55'
$ws.Range("E2").Value = '  +1.36%  '

$ws.Range("D3").Value = '''2.246.53'

$ws.Range("E4").Value = '  +0.05%  '

$ws.Range("D5").Value = '''322.97'
$ws.Range("E5").Value = '  +1.22%  '

$ws.Range("D6").Value = '''101.58'
$ws.Range("E6").Value = '  +1.25%  '

$ws.Range("E7").Value = '  -0.82%  '

$ws.Range("E8").Value = '  +0.07%  '

$ws.Range("E9").Value = '  -1.13%  '

$ws.Range("D10").Value = '''37.15'
$ws.Range("E10").Value = '  +0.71%  '

$ws.Range("D11").Value = '''0.0833'
$ws.Range("E11").Value = '  +0.53%  '

$ws.Range("D12").Value = '''7.73'
$ws.Range("E12").Value = '  +0.47%  '

$ws.Range("E13").Value = '  -2.40%  '

$ws.Range("D14").Value = '''2.588.18'
$ws.Range("E14").Value = '  +0.43%  '

$ws.Range("E15").Value = '  -0.80%  '

$ws.Range("B16").Value = 'WrappedEther'
$ws.Range("C16").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D16").Value = '''2.285.36'
$ws.Range("E16").Value = '  +1.94%  '

$ws.Range("B17").Value = 'Chainlink'
$ws.Range("C17").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D17").Value = '''14.18'
$ws.Range("E17").Value = '  -1.52%  '

$ws.Range("D18").Value = '''43.682.73'
$ws.Range("E18").Value = '  +1.29%  '

$ws.Range("D19").Value = '''13.59'
$ws.Range("E19").Value = '  -6.90%  '

$ws.Range("E20").Value = '  +2.30%  '

$ws.Range("D21").Value = '''6.55'
$ws.Range("E21").Value = '  +0.50%  '

$ws.Range("E22").Value = '  -0.08%  '

$ws.Range("E23").Value = '  -0.77%  '

$ws.Range("D24").Value = '''236.37'
$ws.Range("E24").Value = '  -0.30%  '

$ws.Range("D25").Value = '''2.14'
$ws.Range("E25").Value = '  -0.47%  '

$ws.Range("E26").Value = '  +0.17%  '

$ws.Range("D27").Value = '''10.15'
$ws.Range("E27").Value = '  +1.24%  '

$ws.Range("E28").Value = '  -2.02%  '

$ws.Range("D29").Value = '''37.14'
$ws.Range("E29").Value = '  +4.92%  '

$ws.Range("E30").Value = '  -0.74%  '

$ws.Range("D31").Value = '''161.06'
$ws.Range("E31").Value = '  +5.19%  '

$ws.Range("E32").Value = '  -1.42%  '

$ws.Range("D33").Value = '''0.0851'
$ws.Range("E33").Value = '  -2.39%  '

$ws.Range("E34").Value = '  -1.16%  '

$ws.Range("E35").Value = '  +10.84%  '

$ws.Range("D36").Value = '''3.06'
$ws.Range("E36").Value = '  -3.55%  '

$ws.Range("D37").Value = '''1.93'
$ws.Range("E37").Value = '  -1.96%  '

$ws.Range("E38").Value = '  -2.58%  '

$ws.Range("D39").Value = '''3.76'
$ws.Range("E39").Value = '  +2.68%  '

$ws.Range("D40").Value = '''16.09'
$ws.Range("E40").Value = '  +23.98%  '

$ws.Range("E41").Value = '  -4.31%  '

$ws.Range("E42").Value = '  -1.69%  '

$ws.Range("E43").Value = '  +0.25%  '

$ws.Range("D44").Value = '''1.815.46'
$ws.Range("E44").Value = '  +1.77%  '

$ws.Range("D45").Value = '''76.57'
$ws.Range("E45").Value = '  +1.21%  '

$ws.Range("E46").Value = '  -2.69%  '

$ws.Range("D47").Value = '''82.49'
$ws.Range("E47").Value = '  -4.16%  '

$ws.Range("E48").Value = '  -2.10%  '

$ws.Range("D49").Value = '''58.64'

$ws.Range("D50").Value = '''1.69'
$ws.Range("E50").Value = '  +6.53%  '

$ws.Range("D51").Value = '''103.34'
$ws.Range("E51").Value = '  -0.05%  '
